$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value looks like a plain number need to be forced
# to remain text (matching the source inlineStr cells), without leaving a
# residual custom style behind.
function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = '57.194.53'
$ws.Range("E2").Value = '  +8.89%  '
$ws.Range("D3").Value = '3.255.92'
$ws.Range("E3").Value = '  +4.50%  '
$ws.Range("E4").Value = '  +0.03%  '
Set-TextValue "D5" '394.76'
$ws.Range("E5").Value = '  -0.29%  '
Set-TextValue "D6" '108.55'
$ws.Range("E6").Value = '  +4.46%  '
$ws.Range("D7").Value = '3.250.48'
$ws.Range("E7").Value = '  +4.40%  '
Set-TextValue "D8" '0.566'
$ws.Range("E8").Value = '  +4.70%  '
Set-TextValue "D9" '0.999'
$ws.Range("E9").Value = '  -0.10%  '
Set-TextValue "D10" '0.621'
$ws.Range("E10").Value = '  +3.03%  '
Set-TextValue "D11" '39.23'
$ws.Range("E11").Value = '  +2.78%  '
Set-TextValue "D12" '0.0978'
$ws.Range("E12").Value = '  +13.29%  '
$ws.Range("E13").Value = '  +2.07%  '
$ws.Range("D14").Value = '3.772.80'
$ws.Range("E14").Value = '  +4.35%  '
Set-TextValue "D15" '8.16'
$ws.Range("E15").Value = '  +3.81%  '
Set-TextValue "D16" '19.06'
$ws.Range("E16").Value = '  +1.64%  '
$ws.Range("D17").Value = '3.261.93'
$ws.Range("E17").Value = '  +5.45%  '
$ws.Range("E18").Value = '  -1.25%  '
$ws.Range("E19").Value = '  -2.92%  '
$ws.Range("D20").Value = '57.028.47'
$ws.Range("E20").Value = '  +8.86%  '
$ws.Range("E21").Value = '  +2.58%  '
Set-TextValue "D22" '0.0000106'
$ws.Range("E22").Value = '  +9.14%  '
Set-TextValue "D23" '13.09'
$ws.Range("E23").Value = '  +2.66%  '
Set-TextValue "D24" '300.49'
$ws.Range("E24").Value = '  +11.80%  '
Set-TextValue "D25" '74.18'
$ws.Range("E25").Value = '  +4.62%  '
$ws.Range("E26").Value = '  -2.67%  '
Set-TextValue "D27" '28.03'
$ws.Range("E27").Value = '  +1.84%  '
$ws.Range("B28").Value = 'Filecoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue "D28" '7.92'
$ws.Range("E28").Value = '  -1.75%  '
$ws.Range("B29").Value = 'LEO'
$ws.Range("C29").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue "D29" '4.39'
$ws.Range("E29").Value = '  +3.23%  '
$ws.Range("E30").Value = '  +1.49%  '
$ws.Range("E31").Value = '  -2.41%  '
$ws.Range("E32").Value = '  -0.05%  '
Set-TextValue "D33" '0.109'
$ws.Range("E33").Value = '  +1.87%  '
Set-TextValue "D34" '11.00'
$ws.Range("E34").Value = '  +1.14%  '
Set-TextValue "D35" '37.92'
$ws.Range("E35").Value = '  +3.71%  '
Set-TextValue "D36" '0.0486'
$ws.Range("E36").Value = '  -1.59%  '
$ws.Range("E37").Value = '  +1.83%  '
Set-TextValue "D38" '51.68'
$ws.Range("E38").Value = '  +3.55%  '
$ws.Range("B39").Value = 'FirstDigitalUSD'
$ws.Range("C39").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue "D39" '0.999'
$ws.Range("E39").Value = '  -0.19%  '
$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue "D40" '3.08'
$ws.Range("E40").Value = '  +14.30%  '
Set-TextValue "D41" '3.48'
$ws.Range("E41").Value = '  +2.06%  '
Set-TextValue "D42" '134.14'
$ws.Range("E42").Value = '  +2.96%  '
$ws.Range("E43").Value = '  +2.24%  '
Set-TextValue "D44" '17.21'
$ws.Range("E44").Value = '  +1.48%  '
$ws.Range("E45").Value = '  +3.17%  '
$ws.Range("E46").Value = '  -3.01%  '
Set-TextValue "D47" '0.286'
$ws.Range("E47").Value = '  -2.40%  '
Set-TextValue "D48" '22.01'
$ws.Range("E48").Value = '  -1.18%  '
$ws.Range("D49").Value = '2.147.61'
$ws.Range("E49").Value = '  +2.99%  '
$ws.Range("E50").Value = '  +0.63%  '
$ws.Range("E51").Value = '  -2.71%  '
